# Update "想去人数" (number of people interested) figures in column F
# on both the "展览" (Exhibition) sheet and the "全部类型" (All types) sheet.

$wb = $excel.ActiveWorkbook

# Map of cell address -> new value to apply
$updates = @{
    "F2"  = 1563
    "F4"  = 1020
    "F5"  = 20
    "F7"  = 2601
    "F9"  = 1634
    "F12" = 527
    "F14" = 6
    "F15" = 52
    "F16" = 77
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($addr in $updates.Keys) {
        $ws.Range($addr).Value = $updates[$addr]
    }
}
